$wb = $excel.ActiveWorkbook

# --- "About" sheet: add a new B1 label "California" ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value = "California"

# --- "BECbIC" sheet: refresh the Employee Compensation row with updated RMI data ---
$ws = $wb.Worksheets.Item("BECbIC")
$ws.Range("B2").Value = 1500177304.000151
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 789160759.9395835
$ws.Range("E2").Value = 513924976.8438653
$ws.Range("F2").Value = 642714387.2891195
$ws.Range("G2").Value = 13528928636.45347
$ws.Range("H2").Value = 1812330438.944068
$ws.Range("I2").Value = 1273270653.616102
$ws.Range("J2").Value = 3884649655.032723
$ws.Range("L2").Value = 2700488397.87426
$ws.Range("M2").Value = 7454527965.904436
$ws.Range("N2").Value = 3298801568.164527
$ws.Range("O2").Value = 536565564.8691985
$ws.Range("P2").Value = 1868921075.038227
$ws.Range("Q2").Value = 372512788.5363376
$ws.Range("R2").Value = 407890512.9874496
$ws.Range("S2").Value = 8996406421.030363
$ws.Range("T2").Value = 19917899620.84661
$ws.Range("U2").Value = 3144676318.918274
$ws.Range("V2").Value = 6778381153.518147
$ws.Range("W2").Value = 4034834950.822632
$ws.Range("X2").Value = 8318610437.780892
$ws.Range("Y2").Value = 9951716932.225187
$ws.Range("Z2").Value = 7882846993.049463
$ws.Range("AA2").Value = 2127134217.423578
$ws.Range("AB2").Value = 6617165287.721844
$ws.Range("AC2").Value = 66757363978.61462
$ws.Range("AD2").Value = 87693638645.21896
$ws.Range("AF2").Value = 50681844680.18048
$ws.Range("AG2").Value = 75379993648.03831
$ws.Range("AK2").Value = 23638872077.14451
$ws.Range("AL2").Value = 225677395739.6388
$ws.Range("AP2").Value = 22355885760.72787
$ws.Range("AQ2").Value = 0
